$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"3.090860900556436e-13"
$ws.Range("C2").Value = [double]"2.052465086777033e-09"
$ws.Range("D2").Value = [double]"3.223369029078222"
$ws.Range("E2").Value = [double]"13.86384647080068"
$ws.Range("F2").Value = [double]"0"
$ws.Range("G2").Value = [double]"17.08721550193168"
